$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceAll)
}

# 1. "Nogle skal planlægge skemaet, som har krav til skemaet og dets indhold og andre skal følge det."
#    -> "Nogle skal planlægge skemaet og andre skal følge det."
Replace-Text " skemaet, som har krav til skemaet og dets indhold og andre skal følge det." " skemaet og andre skal følge det."

# 2. Kommunernes mål paragraph
Replace-Text "Både uddannelsesministeriet og kommunen visse krav til hvilke fag der skal skrives på skemaet og hvor mange lektioner der skal afsættes til de forskellige fag. De vil ikke kunne mærke en forskel hvis der kom en software løsning." "Både uddannelsesministeriet og kommunen har krav til hvilke fag der skal skrives på skemaet og hvor mange lektioner der skal afsættes til fagene. De vil ikke kunne mærke en forskel hvis der kom en software løsning, da kommunen ikke planlægger skemaerne."

# 3. "Skolelederen arbejde ud for et budget" -> "Skolelederen arbejder ud for et budget"
Replace-Text "Skolelederen arbejde ud for et budget" "Skolelederen arbejder ud for et budget"

# 4. Derfor ville en softwareløsning ...
Replace-Text ". Derfor ville en softwareløsning som kunne reducere timeantallet det tager at ligge skema have stor interesse hos skolelederen." ". Derfor ville en softwareløsning, som kunne reducere timeantallet det tager at lægge skoleskemaet, have stor interesse hos skolelederen."

# 5. penge andre steder ... Skolelederen har stor indflydelse ...
Replace-Text " penge andre steder på skolen hvor de ville have mere gavn. Skolelederen har stor indflydelse på om programmet nogensinde bliver til noget, da det er skolelederen" " penge andre steder på skolen hvor de ville have mere gavn. Derudover er det skolelederen der bliver kontaktet af lærer, elever eller forældre, hvis skemaet ikke er planlagt korrekt. Skolelederen har stor indflydelse på om programmet bliver en realitet, da det er skolelederen"

# 6. "Så skolelederen er interesseret" -> "Skolelederen er interesseret"
Replace-Text "det løser. Så skolelederen er interesseret" "det løser. Skolelederen er interesseret"

# 7. Eleverne paragraph
Replace-Text "den endelige løsning kommer til at se ud, da de ikke har en fingre med i spillet når det kommer til skemaplanlægning." "den endelige løsning kommer til at se ud, da de ikke er aktive i skemalægningsprocessen."
Replace-Text "og et dårligt lagt skema vil gøre at de F.eks." "og et dårligt planlagt skema vil gøre, at de F.eks."

# 8. Forældrene paragraph
Replace-Text "Forældrene har dog ingen påvirkning på hvordan skemaet bliver lagt, og derfor ikke vil opdage hvis skolen begynder at bruge en softwareløsning til at gøre arbejdet." "Forældrene har dog ingen påvirkning på hvordan skemaet bliver lagt, derfor vil en softwareløsningen ingen effekt have på dem."

# 9. Table cell: Eleverne bliver påvirket meget af skemaet...
#    (split into two replacements so the _GoBack bookmark in the middle of the
#    original text is preserved, ending up after the new text, like in the target)
Replace-Text "Eleverne bliver påv" "Eleverne bliver påvirket meget af skemaet. De har dog ingen indflydelse på hvorledes skemaet bliver planlagt"
Replace-Text "irket meget af skemaet, men har meget lidt påvirkning." ""

# 10. Table cell: "Minimal påvirkning" -> "Lille påvirkning"
Replace-Text "Minimal påvirkning" "Lille påvirkning"

Write-Output "All replacements executed"
